$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Extend the table with a new "2023" column (T), mirroring column S ---
# Copy column S's cell formatting into column T first (so the new cells pick
# up the same visual style as the corresponding S cell), then overwrite the
# values with the new 2023 figures.

$ws.Range("S4").Copy($ws.Range("T4"))
$ws.Range("T4").Value = 2023

$ws.Range("S5").Copy($ws.Range("T5"))
$ws.Range("T5").Value = 43.1

$ws.Range("S6").Copy($ws.Range("T6"))
$ws.Range("T6").Value = 19.7

$ws.Range("S7").Copy($ws.Range("T7"))
$ws.Range("T7").Value = 7.8

$ws.Range("S8").Copy($ws.Range("T8"))
$ws.Range("T8").Value = 15.6

# --- Row height tweaks ---
$ws.Rows("1").RowHeight = 57
$ws.Rows("4").RowHeight = 16.5

# --- Reset the saved view selection back to the default top-left cell ---
$ws.Range("A1").Select() | Out-Null
